$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that must remain TEXT (not auto-converted to number/percentage)
# need NumberFormat forced to Text ("@") before assignment, then the
# transient style reset back to "Normal" so no stray "s" attribute is left.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '330.60'
Set-TextValue $ws.Range('E2') '2.98%'
Set-TextValue $ws.Range('G2') '3'

Set-TextValue $ws.Range('D3') '41.14'
Set-TextValue $ws.Range('E3') '5.04%'
Set-TextValue $ws.Range('G3') '3'

Set-TextValue $ws.Range('D4') '5.699'
Set-TextValue $ws.Range('E4') '-3.75%'
Set-TextValue $ws.Range('G4') '3'

Set-TextValue $ws.Range('D5') '0.08178'
Set-TextValue $ws.Range('E5') '2.13%'
Set-TextValue $ws.Range('G5') '3'

Set-TextValue $ws.Range('D6') '2.074'
Set-TextValue $ws.Range('E6') '9.90%'
Set-TextValue $ws.Range('G6') '3'

$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D7') '4.544'
Set-TextValue $ws.Range('E7') '-0.90%'
Set-TextValue $ws.Range('G7') '3'

$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws.Range('D8') '8.725'
Set-TextValue $ws.Range('E8') '1.05%'
Set-TextValue $ws.Range('G8') '3'

$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D9') '2.949'
Set-TextValue $ws.Range('E9') '0.26%'
Set-TextValue $ws.Range('G9') '3'

$ws.Range('B10').Value = 'MXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D10') '0.9251'
Set-TextValue $ws.Range('E10') '-0.94%'
Set-TextValue $ws.Range('G10') '3'

$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D11') '0.1250'
Set-TextValue $ws.Range('E11') '0.31%'
Set-TextValue $ws.Range('G11') '3'

$ws.Range('B12').Value = 'WazirX'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D12') '0.1957'
Set-TextValue $ws.Range('E12') '0.35%'
Set-TextValue $ws.Range('G12') '3'

$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D13') '0.09433'
Set-TextValue $ws.Range('E13') '2.32%'
Set-TextValue $ws.Range('G13') '3'

$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D14') '0.03691'
Set-TextValue $ws.Range('E14') '5.04%'
Set-TextValue $ws.Range('G14') '3'

$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D15') '0.1055'
Set-TextValue $ws.Range('E15') '10.23%'
Set-TextValue $ws.Range('G15') '3'

$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D16') '0.001303'
Set-TextValue $ws.Range('E16') '1.35%'
Set-TextValue $ws.Range('G16') '3'

$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D17') '0.006251'
Set-TextValue $ws.Range('E17') '2.23%'
Set-TextValue $ws.Range('G17') '3'

$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D18') '3.412'
Set-TextValue $ws.Range('E18') '1.69%'
Set-TextValue $ws.Range('G18') '3'

Set-TextValue $ws.Range('E19') '-1.42%'
Set-TextValue $ws.Range('G19') '3'

Set-TextValue $ws.Range('D20') '8.322'
Set-TextValue $ws.Range('E20') '-4.77%'
Set-TextValue $ws.Range('G20') '3'

Set-TextValue $ws.Range('D21') '0.1380'
Set-TextValue $ws.Range('E21') '-2.80%'
Set-TextValue $ws.Range('G21') '3'

Set-TextValue $ws.Range('D22') '0.2651'
Set-TextValue $ws.Range('E22') '9.85%'
Set-TextValue $ws.Range('G22') '3'

Set-TextValue $ws.Range('D23') '0.04446'
Set-TextValue $ws.Range('E23') '-0.40%'
Set-TextValue $ws.Range('G23') '3'

Set-TextValue $ws.Range('E24') '0.63%'
Set-TextValue $ws.Range('G24') '3'

Set-TextValue $ws.Range('D25') '0.004307'
Set-TextValue $ws.Range('E25') '-1.23%'
Set-TextValue $ws.Range('G25') '3'

Set-TextValue $ws.Range('D26') '0.0001202'
Set-TextValue $ws.Range('E26') '5.17%'
Set-TextValue $ws.Range('G26') '3'

Set-TextValue $ws.Range('G27') '3'

Set-TextValue $ws.Range('G28') '3'

Set-TextValue $ws.Range('G29') '3'

Set-TextValue $ws.Range('G30') '3'

Set-TextValue $ws.Range('G31') '3'

Set-TextValue $ws.Range('G32') '3'

Set-TextValue $ws.Range('G33') '3'

Set-TextValue $ws.Range('G34') '3'

Set-TextValue $ws.Range('G35') '3'

Set-TextValue $ws.Range('G36') '3'

Set-TextValue $ws.Range('G37') '3'

Set-TextValue $ws.Range('G38') '3'

Set-TextValue $ws.Range('D39') '0.02807'
Set-TextValue $ws.Range('E39') '16.49%'
Set-TextValue $ws.Range('G39') '3'

Set-TextValue $ws.Range('D40') '0.05466'
Set-TextValue $ws.Range('E40') '5.69%'
Set-TextValue $ws.Range('G40') '3'

Set-TextValue $ws.Range('D41') '0.007684'
Set-TextValue $ws.Range('E41') '3.39%'
Set-TextValue $ws.Range('G41') '3'

Set-TextValue $ws.Range('D42') '0.009459'
Set-TextValue $ws.Range('E42') '3.82%'
Set-TextValue $ws.Range('G42') '3'

Set-TextValue $ws.Range('D43') '0.1419'
Set-TextValue $ws.Range('E43') '0.98%'
Set-TextValue $ws.Range('G43') '3'

Set-TextValue $ws.Range('D44') '0.002133'
Set-TextValue $ws.Range('E44') '1.77%'
Set-TextValue $ws.Range('G44') '3'

Set-TextValue $ws.Range('D45') '0.01193'
Set-TextValue $ws.Range('E45') '4.33%'
Set-TextValue $ws.Range('G45') '3'

Set-TextValue $ws.Range('D46') '0.00006859'
Set-TextValue $ws.Range('E46') '1.43%'
Set-TextValue $ws.Range('G46') '3'

Set-TextValue $ws.Range('E47') '-0.09%'
Set-TextValue $ws.Range('G47') '3'

Set-TextValue $ws.Range('D48') '0.002283'
Set-TextValue $ws.Range('E48') '60.46%'
Set-TextValue $ws.Range('G48') '3'

Set-TextValue $ws.Range('D49') '0.003234'
Set-TextValue $ws.Range('E49') '7.37%'
Set-TextValue $ws.Range('G49') '3'

Set-TextValue $ws.Range('D50') '0.00002103'
Set-TextValue $ws.Range('E50') '-0.09%'
Set-TextValue $ws.Range('G50') '3'

Set-TextValue $ws.Range('D51') '0.0002003'
Set-TextValue $ws.Range('E51') '-0.09%'
Set-TextValue $ws.Range('G51') '3'
